$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Modificare alocare tema mti2": the student previously assigned to theme
# row 8 (D8, "Pintea Catalina") is reassigned to theme row 14 (D14).
$ws.Range("D8").Value = ""
$ws.Range("D14").Value = "Pintea Cătălina"

# Reflect the new point of interest / scroll position in the sheet view.
$ws.Activate()
$ws.Range("D8").Select()
